$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking price
# strings (e.g. "571.19", "1.00") are written verbatim as text instead
# of being auto-coerced to numbers by Excel's input parser.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.516.26'
$ws.Range("E2").Value = '  -3.23%  '
$ws.Range("D3").Value = '2.599.57'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '571.19'
$ws.Range("E5").Value = '  -4.29%  '
$ws.Range("D6").Value = '154.20'
$ws.Range("E6").Value = '  -2.89%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '0.618'
$ws.Range("E8").Value = '  -3.14%  '
$ws.Range("D9").Value = '2.596.60'
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("E10").Value = '  -7.89%  '
$ws.Range("D11").Value = '5.78'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = '0.376'
$ws.Range("E12").Value = '  -5.31%  '
$ws.Range("D13").Value = '0.156'
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("D14").Value = '27.79'
$ws.Range("E14").Value = '  -4.31%  '
$ws.Range("D15").Value = '3.068.65'
$ws.Range("E15").Value = '  -2.13%  '
$ws.Range("D16").Value = '0.0000178'
$ws.Range("E16").Value = '  -8.03%  '
$ws.Range("D17").Value = '63.372.27'
$ws.Range("E17").Value = '  -3.36%  '
$ws.Range("D18").Value = '2.596.55'
$ws.Range("E18").Value = '  -1.62%  '
$ws.Range("D19").Value = '11.88'
$ws.Range("E19").Value = '  -5.11%  '
$ws.Range("D20").Value = '7.44'
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").Value = '4.46'
$ws.Range("E21").Value = '  -6.49%  '
$ws.Range("D22").Value = '338.90'
$ws.Range("E22").Value = '  -4.03%  '
$ws.Range("D24").Value = '67.06'
$ws.Range("E24").Value = '  -3.78%  '
$ws.Range("D25").Value = '1.78'
$ws.Range("E25").Value = '  +0.74%  '
$ws.Range("D26").Value = '0.0000104'
$ws.Range("E26").Value = '  -7.15%  '
$ws.Range("B27").Value = 'Bittensor'
$ws.Range("C27").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D27").Value = '578.23'
$ws.Range("E27").Value = '  +3.20%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '9.04'
$ws.Range("E28").Value = '  -6.40%  '
$ws.Range("D29").Value = '1.54'
$ws.Range("E29").Value = '  -4.73%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '0.159'
$ws.Range("E31").Value = '  -2.47%  '
$ws.Range("D32").Value = '7.73'
$ws.Range("E32").Value = '  -4.15%  '
$ws.Range("D33").Value = '2.03'
$ws.Range("E33").Value = '  -4.76%  '
$ws.Range("D34").Value = '1.70'
$ws.Range("E34").Value = '  -5.71%  '
$ws.Range("D35").Value = '6.52'
$ws.Range("E35").Value = '  -2.20%  '
$ws.Range("D36").Value = '5.31'
$ws.Range("E36").Value = '  -2.88%  '
$ws.Range("D37").Value = '0.398'
$ws.Range("E37").Value = '  -5.42%  '
$ws.Range("D38").Value = '0.997'
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").Value = '19.51'
$ws.Range("E39").Value = '  -4.73%  '
$ws.Range("D40").Value = '153.90'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("E41").Value = '  -5.59%  '
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").Value = '41.49'
$ws.Range("E43").Value = '  -3.20%  '
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").Value = '156.56'
$ws.Range("E45").Value = '  -2.76%  '
$ws.Range("D46").Value = '23.11'
$ws.Range("E46").Value = '  -0.43%  '
$ws.Range("D47").Value = '3.82'
$ws.Range("E47").Value = '  -6.34%  '
$ws.Range("D48").Value = '0.0578'
$ws.Range("E48").Value = '  -6.11%  '
$ws.Range("D49").Value = '0.625'
$ws.Range("E49").Value = '  -2.82%  '
$ws.Range("D50").Value = '0.0992'
$ws.Range("E50").Value = '  -2.33%  '
$ws.Range("D51").Value = '0.0244'
$ws.Range("E51").Value = '  -5.26%  '

# Restore the original (default/General) formatting on column D now
# that the text values are safely stored, so no style/format change
# is left behind in the saved workbook.
$ws.Range("D2:D51").ClearFormats()
